# Leave Card update (4/12/2023 4:43 PM)
# A new leave entry ("SL(11-0-0)" covering "2/23,27,28, 3/1-10/2023") is
# inserted into the leave table on Sheet1, directly above the existing
# row that starts with the 44986 (1/24/2023) SL(19-0-0) entry. Excel
# shifts every following table row down by one and the table (Table1)
# grows by a row to keep the trailing totals row intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank row above row 132 - this is the row that currently
# holds the 1/24/2023 "SL(19-0-0)" entry; everything below shifts down by one.
$ws.Rows.Item(132).Insert()

# The table's defined range does not auto-grow from a plain row insert,
# so extend Table1 to cover the new last row (the totals row moved from
# 178 to 179).
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K179"))

# The freshly inserted row 132 does not inherit the table's row
# formatting automatically - copy it over from row 133 (the row that
# used to be 132 before the insert).
$ws.Range("A133:K133").Copy()
$ws.Range("A132:K132").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the calculated "EARNED " column formula on the new row (wiped
# out by the formats-only paste above) and re-assert it on the new last
# (totals) row so both keep the structured-table reference form.
$ws.Range("G132").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G179").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# Fill in the new row's data: PARTICULARS + REMARKS for the new SL entry.
$ws.Range("B132").Value = "SL(11-0-0)"
$ws.Range("K132").Value = "2/23,27,28, 3/1-10/2023"
$ws.Range("K132").Style = "Normal"
$ws.Cells.Item(132, 11).Font.Bold = $true

# Update the view so the bottom pane follows the newly inserted row.
$ws.Range("E134").Select()

Write-Output "Leave card row inserted."
